$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the last four "Excel SecurityN" project entries to the new project names.
$ws.Range("A2").Value = "May 23 Project"
$ws.Range("A3").Value = "Tomorrow's Project"
$ws.Range("A4").Value = "Yesterday's Project"
$ws.Range("A5").Value = "May 23 Project"

# Move the active selection from C2:C5 down to A6.
[void]$ws.Range("A6").Select()
